$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$colB = @(1.422295579437503, 1.280489998532573, 1.193256204041575, 1.157668829513455, 1.151757297730171, 1.192776414427726, 1.373436622944723, 1.726312984291724, 1.984621953072974, 2.101908566657414, 2.146288295742295, 2.136731887915118, 2.105560413296303, 2.086462458482117, 1.976952357702089, 1.90971339029511, 1.871018815327716, 1.857914056300274, 1.91687323375811, 2.114717178957505, 2.243819292480339, 2.174934253135689, 1.913636386687131, 1.631010356328431)
$colC = @(0.3967928226192043, 0.3888936859754182, 0.3841169021412867, 0.3821890224420628, 0.3818700365999348, 0.3840908260226712, 0.3940541364432306, 0.4141636668171884, 0.4292746632407898, 0.436219770324243, 0.4388596896160948, 0.4382906969373153, 0.4364367599616799, 0.4353024592791996, 0.4288221915624604, 0.4248647677327426, 0.4225952638943511, 0.4218280078119676, 0.4252853505322207, 0.4369810379684509, 0.4446827716124346, 0.4405669952951996, 0.4250951873622739, 0.4086635647863375)
$colE = @(0.2814513378807391, 0.2826141269567977, 0.28344269311944, 0.2838091995990553, 0.2838718019027517, 0.2834475190519612, 0.2818285040103063, 0.2795613738048885, 0.2784471164613294, 0.2780595370385726, 0.2779298900495988, 0.2779570509402376, 0.278048527971599, 0.2781067889201658, 0.2784748433936954, 0.2787311643774473, 0.2788898247446667, 0.2789454744608015, 0.2787027164223943, 0.2780211945682041, 0.277675558491552, 0.2778509126946354, 0.2787155425481203, 0.2800777169430511)
$colF = @(1.588584802765887, 1.595563293899247, 1.600867218179907, 1.603284807514619, 1.603701718158241, 1.600898785389738, 1.590779430509293, 1.579025925325539, 1.575331065791843, 1.574724845899411, 1.574649925714695, 1.574659181524723, 1.574715582248686, 1.574770271731424, 1.575392317280361, 1.576049234720415, 1.576528214917829, 1.57670775643048, 1.575968836406759, 1.574694818123035, 1.574763598916974, 1.574644370698266, 1.576004868921331, 1.581338553740906)
$colG = @(0.4192116588298305, 0.4247075550293573, 0.4285082427658367, 0.4301639500946379, 0.4304453273211308, 0.4285301397567665, 0.4210180261315131, 0.4096806908093455, 0.4034372819470988, 0.4010536438018377, 0.4002169638013697, 0.4003942197382102, 0.4009834861172834, 0.4013530267113836, 0.4036022707706977, 0.4050992524748622, 0.4060032284505795, 0.4063166668543232, 0.404935448661611, 0.4008086119565917, 0.3984960447494075, 0.3996950179065237, 0.4050093693834, 0.4123826241993314)
$colH = @(0.5862643707296371, 0.593302144479118, 0.597968859211818, 0.5999574242195322, 0.6002928678345896, 0.5979953261287179, 0.5886192666580143, 0.5729757456222231, 0.5631566822845073, 0.5590537317899944, 0.5575524078653515, 0.5578734144892081, 0.5589291669539023, 0.5595826683409797, 0.5634321454284148, 0.5658868944475728, 0.5673330391786351, 0.5678285574661928, 0.5656220382683586, 0.5586176449468852, 0.5543451754678657, 0.5565975169389077, 0.5657416711114109, 0.5769138698791494)
$colJ = @(0.03502352788252239, 0.03295554047478788, 0.0316850409194771, 0.03116715670934056, 0.03108115490987018, 0.03167805708383398, 0.03431066561184792, 0.0394656699514826, 0.04324668330933434, 0.04496503229143656, 0.04561545150708213, 0.0454753853508123, 0.04501854858251164, 0.04473868479671239, 0.04313434808436512, 0.04214968430657251, 0.04158317880371243, 0.041391344869524, 0.04225451959931803, 0.04515274062909924, 0.04704522998824956, 0.04603533984426633, 0.042207124834583, 0.03807209602600636)
$colL = @(0.5591069227842951, 0.5421218507654828, 0.5318377106014651, 0.5276835916519786, 0.5269960334080253, 0.5317815374281736, 0.5532206278183622, 0.596398344271222, 0.6287996803612543, 0.6436845970411866, 0.6493417351478001, 0.6481224613368113, 0.6441496035551495, 0.6417187780345728, 0.6278298105087003, 0.6193463567434776, 0.6144806045677171, 0.6128355096089706, 0.6202480189917594, 0.6453159735382599, 0.6618189052249761, 0.6530001584716558, 0.6198403416479437, 0.5845974644336849)
$colO = @(1.960301203518554, 1.986406789389221, 2.004052987939659, 2.011649973464628, 2.012935950179028, 2.004153800475422, 1.968966352618608, 1.912825795930118, 1.879461360462557, 1.866003327897502, 1.861155123278877, 1.862188225947492, 1.865599487228963, 1.867721311142873, 1.880375544527496, 1.888579499353227, 1.893460009183926, 1.895140233441282, 1.887689421743488, 1.864590778610847, 1.850940724401241, 1.858093426735067, 1.888091315203198, 1.926632087537783)

$columns = @("B", "C", "E", "F", "G", "H", "J", "L", "O")
$colArrays = @{
    "B" = $colB
    "C" = $colC
    "E" = $colE
    "F" = $colF
    "G" = $colG
    "H" = $colH
    "J" = $colJ
    "L" = $colL
    "O" = $colO
}

foreach ($col in $columns) {
    $vals = $colArrays[$col]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $row = $i + 2
        $ws.Range("$col$row").Value = $vals[$i]
    }
}
